$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 5; existing rows 5-17 shift down to 6-18.
$ws.Rows("5:5").Insert()

# Populate the newly inserted row 5 with the new weekly price record.
$ws.Cells.Item(5,1).Value = 1
$ws.Cells.Item(5,2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(5,3).Value = "Arica y Parinacota"
$ws.Cells.Item(5,4).Value = 44649
$ws.Cells.Item(5,5).Value = 15
$ws.Cells.Item(5,6).Value = 100112026
$ws.Cells.Item(5,7).Value = "Haba"
$ws.Cells.Item(5,8).Value = "Sin especificar"
$ws.Cells.Item(5,9).Value = "Primera"
$ws.Cells.Item(5,10).Value = 600
$ws.Cells.Item(5,11).Value = 900
$ws.Cells.Item(5,12).Value = 1000
$ws.Cells.Item(5,13).Value = 950
$ws.Cells.Item(5,14).Value = "`$/kilo"
$ws.Cells.Item(5,15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(5,16).Value = 950
$ws.Cells.Item(5,17).Value = 1
$ws.Cells.Item(5,18).Value = "Hortaliza"
